# Remove the "Rodriguez Garcia" (row 85) and "Webb" (row 103) author records.
# Delete from the bottom up so row numbers above the deleted row are not
# disturbed before they are processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(103).Delete()
$ws.Rows.Item(85).Delete()
